$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" '243.88'
Set-TextValue "E2" '-0.69%'
Set-TextValue "G2" '8'
Set-TextValue "D3" '26.70'
Set-TextValue "E3" '4.67%'
Set-TextValue "G3" '8'
Set-TextValue "D4" '5.139'
Set-TextValue "E4" '0.08%'
Set-TextValue "G4" '8'
Set-TextValue "D5" '0.05614'
Set-TextValue "G5" '8'
Set-TextValue "D6" '6.466'
Set-TextValue "E6" '-0.38%'
Set-TextValue "G6" '8'
Set-TextValue "D7" '0.8196'
Set-TextValue "E7" '0.13%'
Set-TextValue "G7" '8'
Set-TextValue "D8" '0.8322'
Set-TextValue "E8" '-1.96%'
Set-TextValue "G8" '8'
Set-TextValue "D9" '0.1330'
Set-TextValue "E9" '-0.73%'
Set-TextValue "G9" '8'
Set-TextValue "D10" '0.06932'
Set-TextValue "E10" '-0.48%'
Set-TextValue "G10" '8'
Set-TextValue "D11" '0.02893'
Set-TextValue "E11" '0.83%'
Set-TextValue "G11" '8'
Set-TextValue "D12" '0.09383'
Set-TextValue "E12" '-0.02%'
Set-TextValue "G12" '8'
Set-TextValue "D13" '0.001513'
Set-TextValue "E13" '-0.46%'
Set-TextValue "G13" '8'
Set-TextValue "D14" '0.0005998'
Set-TextValue "E14" '-93.83%'
Set-TextValue "G14" '8'
Set-TextValue "D15" '0.006220'
Set-TextValue "E15" '1.20%'
Set-TextValue "G15" '8'
Set-TextValue "E16" '3.35%'
Set-TextValue "G16" '8'
Set-TextValue "D17" '3.025'
Set-TextValue "E17" '0.05%'
Set-TextValue "G17" '8'
Set-TextValue "D18" '2.301'
Set-TextValue "E18" '7.34%'
Set-TextValue "G18" '8'
Set-TextValue "E19" '-1.68%'
Set-TextValue "G19" '8'
Set-TextValue "D20" '0.03088'
Set-TextValue "E20" '-4.60%'
Set-TextValue "G20" '8'
Set-TextValue "D21" '0.1291'
Set-TextValue "E21" '-2.14%'
Set-TextValue "G21" '8'
Set-TextValue "D22" '3.751'
Set-TextValue "E22" '0.12%'
Set-TextValue "G22" '8'
Set-TextValue "D23" '0.04595'
Set-TextValue "E23" '-2.21%'
Set-TextValue "G23" '8'
Set-TextValue "D24" '0.1341'
Set-TextValue "E24" '-2.48%'
Set-TextValue "G24" '8'
Set-TextValue "D25" '0.001226'
Set-TextValue "E25" '-1.66%'
Set-TextValue "G25" '8'
Set-TextValue "D26" '0.004494'
Set-TextValue "E26" '-2.41%'
Set-TextValue "G26" '8'
Set-TextValue "D27" '0.00009601'
Set-TextValue "E27" '0.05%'
Set-TextValue "G27" '8'
Set-TextValue "E28" '0.67%'
Set-TextValue "G28" '8'
Set-TextValue "G29" '8'
Set-TextValue "G30" '8'
Set-TextValue "G31" '8'
Set-TextValue "G32" '8'
Set-TextValue "G33" '8'
Set-TextValue "G34" '8'
Set-TextValue "G35" '8'
Set-TextValue "G36" '8'
Set-TextValue "G37" '8'
Set-TextValue "G38" '8'
Set-TextValue "G39" '8'
Set-TextValue "D40" '0.03637'
Set-TextValue "E40" '-0.42%'
Set-TextValue "G40" '8'
Set-TextValue "B41" 'KickToken'
Set-TextValue "C41" 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue "D41" '0.006171'
Set-TextValue "E41" '1.03%'
Set-TextValue "G41" '8'
Set-TextValue "B42" 'BKEXToken'
Set-TextValue "C42" 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D42" '0.1051'
Set-TextValue "E42" '-0.10%'
Set-TextValue "G42" '8'
Set-TextValue "B43" 'CEJI'
Set-TextValue "C43" 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue "D43" '0.002500'
Set-TextValue "E43" '0.05%'
Set-TextValue "G43" '8'
Set-TextValue "D44" '0.008112'
Set-TextValue "E44" '4.29%'
Set-TextValue "G44" '8'
Set-TextValue "D45" '0.00005346'
Set-TextValue "E45" '0.58%'
Set-TextValue "G45" '8'
Set-TextValue "E46" '0.05%'
Set-TextValue "G46" '8'
Set-TextValue "E47" '-18.35%'
Set-TextValue "G47" '8'
Set-TextValue "D48" '0.002550'
Set-TextValue "E48" '19.95%'
Set-TextValue "G48" '8'
Set-TextValue "E49" '0.05%'
Set-TextValue "G49" '8'
Set-TextValue "E50" '0.05%'
Set-TextValue "G50" '8'
Set-TextValue "G51" '8'
